$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 5.336433666666667
$ws.Cells.Item(2, 8).Value2 = 16.009301
$ws.Cells.Item(2, 9).Value2 = 0.09554123656860429
$ws.Cells.Item(2, 10).Value2 = 0.09601942232049432
$ws.Cells.Item(2, 13).Value2 = 38.940166
$ws.Cells.Item(2, 14).Value2 = 116.820498
$ws.Cells.Item(2, 15).Value2 = 0.1629440517615001
$ws.Cells.Item(2, 16).Value2 = 0.1644242630434466
$ws.Cells.Item(2, 17).Value2 = 207.8016128279887
$ws.Cells.Item(2, 18).Value2 = 1870.214515451898
$ws.Cells.Item(2, 19).Value2 = 0.01556787619679239
$ws.Cells.Item(2, 20).Value2 = 0.01578792275290474
$ws.Cells.Item(3, 7).Value2 = 5.336433666666667
$ws.Cells.Item(3, 8).Value2 = 16.009301
$ws.Cells.Item(3, 9).Value2 = 0.09554123656860429
$ws.Cells.Item(3, 10).Value2 = 0.09601942232049432
$ws.Cells.Item(3, 14).Value2 = 7.785288
$ws.Cells.Item(3, 15).Value2 = 0.01085910771284493
$ws.Cells.Item(3, 16).Value2 = 0.01095775368104481
$ws.Cells.Item(3, 17).Value2 = 13.848557662632
$ws.Cells.Item(3, 18).Value2 = 124.637018963688
$ws.Cells.Item(3, 19).Value2 = 0.001037492578916873
$ws.Cells.Item(3, 20).Value2 = 0.001052157178384193
$ws.Cells.Item(4, 7).Value2 = 5.336433666666667
$ws.Cells.Item(4, 8).Value2 = 16.009301
$ws.Cells.Item(4, 9).Value2 = 0.09554123656860429
$ws.Cells.Item(4, 10).Value2 = 0.09601942232049432
$ws.Cells.Item(4, 13).Value2 = 115.357885
$ws.Cells.Item(4, 14).Value2 = 346.073655
$ws.Cells.Item(4, 15).Value2 = 0.4827118914833898
$ws.Cells.Item(4, 16).Value2 = 0.487096927819354
$ws.Cells.Item(4, 17).Value2 = 615.5997012294617
$ws.Cells.Item(4, 18).Value2 = 5540.397311065155
$ws.Cells.Item(4, 19).Value2 = 0.04611889101869299
$ws.Cells.Item(4, 20).Value2 = 0.04677076562330189
$ws.Cells.Item(5, 7).Value2 = 5.336433666666667
$ws.Cells.Item(5, 8).Value2 = 16.009301
$ws.Cells.Item(5, 9).Value2 = 0.09554123656860429
$ws.Cells.Item(5, 10).Value2 = 0.09601942232049432
$ws.Cells.Item(5, 13).Value2 = 6.4541395
$ws.Cells.Item(5, 14).Value2 = 12.908279
$ws.Cells.Item(5, 15).Value2 = 0.02700716891561127
$ws.Cells.Item(5, 16).Value2 = 0.01816833773242602
$ws.Cells.Item(5, 17).Value2 = 34.44208731716317
$ws.Cells.Item(5, 18).Value2 = 206.652523902979
$ws.Cells.Item(5, 19).Value2 = 0.002580298314414673
$ws.Cells.Item(5, 20).Value2 = 0.001744513293591186
$ws.Cells.Item(6, 7).Value2 = 5.336433666666667
$ws.Cells.Item(6, 8).Value2 = 16.009301
$ws.Cells.Item(6, 9).Value2 = 0.09554123656860429
$ws.Cells.Item(6, 10).Value2 = 0.09601942232049432
$ws.Cells.Item(6, 13).Value2 = 75.63146466666666
$ws.Cells.Item(6, 14).Value2 = 226.894394
$ws.Cells.Item(6, 15).Value2 = 0.3164777801266539
$ws.Cells.Item(6, 16).Value2 = 0.3193527177237286
$ws.Cells.Item(6, 17).Value2 = 403.6022943065104
$ws.Cells.Item(6, 18).Value2 = 3632.420648758594
$ws.Cells.Item(6, 19).Value2 = 0.03023667845978737
$ws.Cells.Item(6, 20).Value2 = 0.03066406347231231
$ws.Cells.Item(7, 9).Value2 = 0.0653459693259494
$ws.Cells.Item(7, 10).Value2 = 0.06567302717654226
$ws.Cells.Item(7, 13).Value2 = 38.940166
$ws.Cells.Item(7, 14).Value2 = 116.820498
$ws.Cells.Item(7, 15).Value2 = 0.1629440517615001
$ws.Cells.Item(7, 16).Value2 = 0.1644242630434466
$ws.Cells.Item(7, 17).Value2 = 142.127088840744
$ws.Cells.Item(7, 18).Value2 = 1279.143799566696
$ws.Cells.Item(7, 19).Value2 = 0.0106477370082529
$ws.Cells.Item(7, 20).Value2 = 0.0107982390953352
$ws.Cells.Item(8, 9).Value2 = 0.0653459693259494
$ws.Cells.Item(8, 10).Value2 = 0.06567302717654226
$ws.Cells.Item(8, 14).Value2 = 7.785288
$ws.Cells.Item(8, 15).Value2 = 0.01085910771284493
$ws.Cells.Item(8, 16).Value2 = 0.01095775368104481
$ws.Cells.Item(8, 18).Value2 = 85.24619431977601
$ws.Cells.Item(8, 19).Value2 = 0.0007095989195107455
$ws.Cells.Item(8, 20).Value2 = 0.0007196288552891121
$ws.Cells.Item(9, 9).Value2 = 0.0653459693259494
$ws.Cells.Item(9, 10).Value2 = 0.06567302717654226
$ws.Cells.Item(9, 13).Value2 = 115.357885
$ws.Cells.Item(9, 14).Value2 = 346.073655
$ws.Cells.Item(9, 15).Value2 = 0.4827118914833898
$ws.Cells.Item(9, 16).Value2 = 0.487096927819354
$ws.Cells.Item(9, 17).Value2 = 421.04289873534
$ws.Cells.Item(9, 18).Value2 = 3789.38608861806
$ws.Cells.Item(9, 19).Value2 = 0.0315432764541446
$ws.Cells.Item(9, 20).Value2 = 0.03198912977829068
$ws.Cells.Item(10, 9).Value2 = 0.0653459693259494
$ws.Cells.Item(10, 10).Value2 = 0.06567302717654226
$ws.Cells.Item(10, 13).Value2 = 6.4541395
$ws.Cells.Item(10, 14).Value2 = 12.908279
$ws.Cells.Item(10, 15).Value2 = 0.02700716891561127
$ws.Cells.Item(10, 16).Value2 = 0.01816833773242602
$ws.Cells.Item(10, 17).Value2 = 23.556860494818
$ws.Cells.Item(10, 18).Value2 = 141.341162968908
$ws.Cells.Item(10, 19).Value2 = 0.001764809631540268
$ws.Cells.Item(10, 20).Value2 = 0.001193169737654212
$ws.Cells.Item(11, 9).Value2 = 0.0653459693259494
$ws.Cells.Item(11, 10).Value2 = 0.06567302717654226
$ws.Cells.Item(11, 13).Value2 = 75.63146466666666
$ws.Cells.Item(11, 14).Value2 = 226.894394
$ws.Cells.Item(11, 15).Value2 = 0.3164777801266539
$ws.Cells.Item(11, 16).Value2 = 0.3193527177237286
$ws.Cells.Item(11, 17).Value2 = 276.046072783432
$ws.Cells.Item(11, 18).Value2 = 2484.414655050888
$ws.Cells.Item(11, 19).Value2 = 0.02068054731250088
$ws.Cells.Item(11, 20).Value2 = 0.02097285970997306
$ws.Cells.Item(12, 7).Value2 = 22.33109633333333
$ws.Cells.Item(12, 8).Value2 = 66.993289
$ws.Cells.Item(12, 9).Value2 = 0.3998064420712607
$ws.Cells.Item(12, 10).Value2 = 0.4018074811092581
$ws.Cells.Item(12, 13).Value2 = 38.940166
$ws.Cells.Item(12, 14).Value2 = 116.820498
$ws.Cells.Item(12, 15).Value2 = 0.1629440517615001
$ws.Cells.Item(12, 16).Value2 = 0.1644242630434466
$ws.Cells.Item(12, 17).Value2 = 869.5765981819914
$ws.Cells.Item(12, 18).Value2 = 7826.189383637922
$ws.Cells.Item(12, 19).Value2 = 0.06514608159144071
$ws.Cells.Item(12, 20).Value2 = 0.06606689896673335
$ws.Cells.Item(13, 7).Value2 = 22.33109633333333
$ws.Cells.Item(13, 8).Value2 = 66.993289
$ws.Cells.Item(13, 9).Value2 = 0.3998064420712607
$ws.Cells.Item(13, 10).Value2 = 0.4018074811092581
$ws.Cells.Item(13, 14).Value2 = 7.785288
$ws.Cells.Item(13, 15).Value2 = 0.01085910771284493
$ws.Cells.Item(13, 16).Value2 = 0.01095775368104481
$ws.Cells.Item(13, 17).Value2 = 57.95133877024801
$ws.Cells.Item(13, 18).Value2 = 521.5620489322321
$ws.Cells.Item(13, 19).Value2 = 0.004341541218741118
$ws.Cells.Item(13, 20).Value2 = 0.004402907405196318
$ws.Cells.Item(14, 7).Value2 = 22.33109633333333
$ws.Cells.Item(14, 8).Value2 = 66.993289
$ws.Cells.Item(14, 9).Value2 = 0.3998064420712607
$ws.Cells.Item(14, 10).Value2 = 0.4018074811092581
$ws.Cells.Item(14, 13).Value2 = 115.357885
$ws.Cells.Item(14, 14).Value2 = 346.073655
$ws.Cells.Item(14, 15).Value2 = 0.4827118914833898
$ws.Cells.Item(14, 16).Value2 = 0.487096927819354
$ws.Cells.Item(14, 17).Value2 = 2576.068042744588
$ws.Cells.Item(14, 18).Value2 = 23184.61238470129
$ws.Cells.Item(14, 19).Value2 = 0.1929913238794626
$ws.Cells.Item(14, 20).Value2 = 0.1957191896231527
$ws.Cells.Item(15, 7).Value2 = 22.33109633333333
$ws.Cells.Item(15, 8).Value2 = 66.993289
$ws.Cells.Item(15, 9).Value2 = 0.3998064420712607
$ws.Cells.Item(15, 10).Value2 = 0.4018074811092581
$ws.Cells.Item(15, 13).Value2 = 6.4541395
$ws.Cells.Item(15, 14).Value2 = 12.908279
$ws.Cells.Item(15, 15).Value2 = 0.02700716891561127
$ws.Cells.Item(15, 16).Value2 = 0.01816833773242602
$ws.Cells.Item(15, 17).Value2 = 144.1280109232719
$ws.Cells.Item(15, 18).Value2 = 864.768065539631
$ws.Cells.Item(15, 19).Value2 = 0.01079764011456809
$ws.Cells.Item(15, 20).Value2 = 0.007300174020208387
$ws.Cells.Item(16, 7).Value2 = 22.33109633333333
$ws.Cells.Item(16, 8).Value2 = 66.993289
$ws.Cells.Item(16, 9).Value2 = 0.3998064420712607
$ws.Cells.Item(16, 10).Value2 = 0.4018074811092581
$ws.Cells.Item(16, 13).Value2 = 75.63146466666666
$ws.Cells.Item(16, 14).Value2 = 226.894394
$ws.Cells.Item(16, 15).Value2 = 0.3164777801266539
$ws.Cells.Item(16, 16).Value2 = 0.3193527177237286
$ws.Cells.Item(16, 17).Value2 = 1688.933523302429
$ws.Cells.Item(16, 18).Value2 = 15200.40170972187
$ws.Cells.Item(16, 19).Value2 = 0.1265298552670482
$ws.Cells.Item(16, 20).Value2 = 0.1283183110939673
$ws.Cells.Item(17, 7).Value2 = 0.8344860000000001
$ws.Cells.Item(17, 8).Value2 = 1.668972
$ws.Cells.Item(17, 9).Value2 = 0.01494028209086487
$ws.Cells.Item(17, 10).Value2 = 0.01001003899602363
$ws.Cells.Item(17, 13).Value2 = 38.940166
$ws.Cells.Item(17, 14).Value2 = 116.820498
$ws.Cells.Item(17, 15).Value2 = 0.1629440517615001
$ws.Cells.Item(17, 16).Value2 = 0.1644242630434466
$ws.Cells.Item(17, 17).Value2 = 32.495023364676
$ws.Cells.Item(17, 18).Value2 = 194.970140188056
$ws.Cells.Item(17, 19).Value2 = 0.002434430098345298
$ws.Cells.Item(17, 20).Value2 = 0.001645893284957347
$ws.Cells.Item(18, 7).Value2 = 0.8344860000000001
$ws.Cells.Item(18, 8).Value2 = 1.668972
$ws.Cells.Item(18, 9).Value2 = 0.01494028209086487
$ws.Cells.Item(18, 10).Value2 = 0.01001003899602363
$ws.Cells.Item(18, 14).Value2 = 7.785288
$ws.Cells.Item(18, 15).Value2 = 0.01085910771284493
$ws.Cells.Item(18, 16).Value2 = 0.01095775368104481
$ws.Cells.Item(18, 17).Value2 = 2.165571280656001
$ws.Cells.Item(18, 18).Value2 = 12.993427683936
$ws.Cells.Item(18, 19).Value2 = 0.0001622381324849897
$ws.Cells.Item(18, 20).Value2 = 0.00010968754165608
$ws.Cells.Item(19, 7).Value2 = 0.8344860000000001
$ws.Cells.Item(19, 8).Value2 = 1.668972
$ws.Cells.Item(19, 9).Value2 = 0.01494028209086487
$ws.Cells.Item(19, 10).Value2 = 0.01001003899602363
$ws.Cells.Item(19, 13).Value2 = 115.357885
$ws.Cells.Item(19, 14).Value2 = 346.073655
$ws.Cells.Item(19, 15).Value2 = 0.4827118914833898
$ws.Cells.Item(19, 16).Value2 = 0.487096927819354
$ws.Cells.Item(19, 17).Value2 = 96.26454002211
$ws.Cells.Item(19, 18).Value2 = 577.58724013266
$ws.Cells.Item(19, 19).Value2 = 0.007211851827376792
$ws.Cells.Item(19, 20).Value2 = 0.004875859242315039
$ws.Cells.Item(20, 7).Value2 = 0.8344860000000001
$ws.Cells.Item(20, 8).Value2 = 1.668972
$ws.Cells.Item(20, 9).Value2 = 0.01494028209086487
$ws.Cells.Item(20, 10).Value2 = 0.01001003899602363
$ws.Cells.Item(20, 13).Value2 = 6.4541395
$ws.Cells.Item(20, 14).Value2 = 12.908279
$ws.Cells.Item(20, 15).Value2 = 0.02700716891561127
$ws.Cells.Item(20, 16).Value2 = 0.01816833773242602
$ws.Cells.Item(20, 17).Value2 = 5.385889054797
$ws.Cells.Item(20, 18).Value2 = 21.543556219188
$ws.Cells.Item(20, 19).Value2 = 0.0004034947220748693
$ws.Cells.Item(20, 20).Value2 = 0.0001818657691945119
$ws.Cells.Item(21, 7).Value2 = 0.8344860000000001
$ws.Cells.Item(21, 8).Value2 = 1.668972
$ws.Cells.Item(21, 9).Value2 = 0.01494028209086487
$ws.Cells.Item(21, 10).Value2 = 0.01001003899602363
$ws.Cells.Item(21, 13).Value2 = 75.63146466666666
$ws.Cells.Item(21, 14).Value2 = 226.894394
$ws.Cells.Item(21, 15).Value2 = 0.3164777801266539
$ws.Cells.Item(21, 16).Value2 = 0.3193527177237286
$ws.Cells.Item(21, 17).Value2 = 63.113398423828
$ws.Cells.Item(21, 18).Value2 = 378.680390542968
$ws.Cells.Item(21, 19).Value2 = 0.004728267310582916
$ws.Cells.Item(21, 20).Value2 = 0.003196733157900649
$ws.Cells.Item(22, 7).Value2 = 23.70286866666666
$ws.Cells.Item(22, 8).Value2 = 71.10860599999999
$ws.Cells.Item(22, 9).Value2 = 0.4243660699433207
$ws.Cells.Item(22, 10).Value2 = 0.4264900303976816
$ws.Cells.Item(22, 13).Value2 = 38.940166
$ws.Cells.Item(22, 14).Value2 = 116.820498
$ws.Cells.Item(22, 15).Value2 = 0.1629440517615001
$ws.Cells.Item(22, 16).Value2 = 0.1644242630434466
$ws.Cells.Item(22, 17).Value2 = 922.9936405561986
$ws.Cells.Item(22, 18).Value2 = 8306.942765005788
$ws.Cells.Item(22, 19).Value2 = 0.06914792686666883
$ws.Cells.Item(22, 20).Value2 = 0.07012530894351594
$ws.Cells.Item(23, 7).Value2 = 23.70286866666666
$ws.Cells.Item(23, 8).Value2 = 71.10860599999999
$ws.Cells.Item(23, 9).Value2 = 0.4243660699433207
$ws.Cells.Item(23, 10).Value2 = 0.4264900303976816
$ws.Cells.Item(23, 14).Value2 = 7.785288
$ws.Cells.Item(23, 15).Value2 = 0.01085910771284493
$ws.Cells.Item(23, 16).Value2 = 0.01095775368104481
$ws.Cells.Item(23, 17).Value2 = 61.511219665392
$ws.Cells.Item(23, 18).Value2 = 553.600976988528
$ws.Cells.Item(23, 19).Value2 = 0.004608236863191206
$ws.Cells.Item(23, 20).Value2 = 0.004673372700519111
$ws.Cells.Item(24, 7).Value2 = 23.70286866666666
$ws.Cells.Item(24, 8).Value2 = 71.10860599999999
$ws.Cells.Item(24, 9).Value2 = 0.4243660699433207
$ws.Cells.Item(24, 10).Value2 = 0.4264900303976816
$ws.Cells.Item(24, 13).Value2 = 115.357885
$ws.Cells.Item(24, 14).Value2 = 346.073655
$ws.Cells.Item(24, 15).Value2 = 0.4827118914833898
$ws.Cells.Item(24, 16).Value2 = 0.487096927819354
$ws.Cells.Item(24, 17).Value2 = 2734.312797819436
$ws.Cells.Item(24, 18).Value2 = 24608.81518037493
$ws.Cells.Item(24, 19).Value2 = 0.2048465483037128
$ws.Cells.Item(24, 20).Value2 = 0.2077419835522936
$ws.Cells.Item(25, 7).Value2 = 23.70286866666666
$ws.Cells.Item(25, 8).Value2 = 71.10860599999999
$ws.Cells.Item(25, 9).Value2 = 0.4243660699433207
$ws.Cells.Item(25, 10).Value2 = 0.4264900303976816
$ws.Cells.Item(25, 13).Value2 = 6.4541395
$ws.Cells.Item(25, 14).Value2 = 12.908279
$ws.Cells.Item(25, 15).Value2 = 0.02700716891561127
$ws.Cells.Item(25, 16).Value2 = 0.01816833773242602
$ws.Cells.Item(25, 17).Value2 = 152.9816209248457
$ws.Cells.Item(25, 18).Value2 = 917.889725549074
$ws.Cells.Item(25, 19).Value2 = 0.01146092613301337
$ws.Cells.Item(25, 20).Value2 = 0.007748614911777718
$ws.Cells.Item(26, 7).Value2 = 23.70286866666666
$ws.Cells.Item(26, 8).Value2 = 71.10860599999999
$ws.Cells.Item(26, 9).Value2 = 0.4243660699433207
$ws.Cells.Item(26, 10).Value2 = 0.4264900303976816
$ws.Cells.Item(26, 13).Value2 = 75.63146466666666
$ws.Cells.Item(26, 14).Value2 = 226.894394
$ws.Cells.Item(26, 15).Value2 = 0.3164777801266539
$ws.Cells.Item(26, 16).Value2 = 0.3193527177237286
$ws.Cells.Item(26, 17).Value2 = 1792.68267406164
$ws.Cells.Item(26, 18).Value2 = 16134.14406655476
$ws.Cells.Item(26, 19).Value2 = 0.1343024317767345
$ws.Cells.Item(26, 20).Value2 = 0.1283183110939673
